$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Multiplicador Imposto") had kept a stale 1.3 value for the
# Serviço rows while the correct 1.5 tax multiplier had been duplicated
# into a separate column F ("Modificador Imposto"). Fix column D with the
# correct multiplier.
$ws.Range("D3").Value = 1.5
$ws.Range("D5").Value = 1.5
$ws.Range("D7").Value = 1.5
$ws.Range("D8").Value = 1.5

# Column E ("Preço Base Reais") was left blank for the Produto rows.
# Fill it in as Preço Base Original (B) * Multiplicador Imposto (D).
$ws.Range("E2").Value = 1099.989
$ws.Range("E4").Value = 989.9890000000001
$ws.Range("E6").Value = 3300

# Column F ("Modificador Imposto") duplicated column D's information and
# is removed entirely.
$ws.Columns.Item(6).Delete()
